# 1st changes of mifos to finflux
#
# On the "Repayment Schedule" sheet a new (empty) column is inserted
# immediately before column N ("Late"), pushing the existing N/O/P
# columns ("Late" and "Outstanding") one column to the right.
# The previously active sheet ("Transactions") is de-activated and the
# "Repayment Schedule" sheet becomes the active tab/sheet, with a new
# cell selection on it.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N, shifting N:P -> O:Q
$wsSchedule.Columns("N").Insert() | Out-Null

# Make "Repayment Schedule" the active sheet (was "Transactions")
$wsSchedule.Activate() | Out-Null

# Update the selection on the now-active sheet
$wsSchedule.Range("T6").Select() | Out-Null
